$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Name) to hold the new "title" header,
# shifting the existing Name/Email columns (and subsequent ones) to the right.
$ws.Range("B1").EntireColumn.Insert()

# Set the new header values for row 1 across A1:H1
$ws.Range("A1").Value = "updated_date"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "company"
$ws.Range("F1").Value = "event"
$ws.Range("G1").Value = "linkedin"
$ws.Range("H1").Value = "website"
